# Adds a new "2022-Q1" sheet (cloned from the "2021-Q4" sheet, whose layout/
# header/style already match the desired new sheet) positioned right before
# the "总计" (totals) sheet, fills in the quarter's fund data, and updates
# the "总计" sheet with a new leading row summarizing 2022-Q1 while shifting
# the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: Insert the new "2022-Q1" worksheet before "总计" by copying the
# "2021-Q4" worksheet (same column layout/headers/styles as the new sheet).
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetRef = $wb.Worksheets.Item("总计")
$template.Copy($totalSheetRef, [System.Reflection.Missing]::Value)

$newSheet = $wb.Worksheets.Item(4)
$newSheet.Name = "2022-Q1"

# The Copy() operation makes the new sheet the active tab; restore the
# workbook's original active tab (the first sheet) so we don't leave a
# stray tabSelected flag on the new sheet.
$wb.Worksheets.Item(1).Activate()

# Fill in the 2022-Q1 fund figures (columns D:H). Columns A:C already match
# because they were copied from the template sheet.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $newSheet.Cells.Item(2, 4) "35.75"
Set-TextValue $newSheet.Cells.Item(2, 5) "86.53"
Set-TextValue $newSheet.Cells.Item(2, 6) "3.94"
Set-TextValue $newSheet.Cells.Item(2, 7) "1.4086"
$newSheet.Cells.Item(2, 8).Value = 4

Set-TextValue $newSheet.Cells.Item(3, 4) "35.75"
Set-TextValue $newSheet.Cells.Item(3, 5) "86.53"
Set-TextValue $newSheet.Cells.Item(3, 6) "3.94"
Set-TextValue $newSheet.Cells.Item(3, 7) "1.4086"
$newSheet.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------
# Step 2: Update the "总计" worksheet - add a new 2022-Q1 row at the top
# of the data (row 2) and push the existing rows down by one.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Create row 5 (new) by copying the formatting of row 4 so the styled "A"
# column cell keeps its formatting.
$totals.Range("A4:D4").Copy()
$totals.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Shift existing rows down (bottom-up to avoid clobbering data we still need).
$totals.Cells.Item(5, 1).Value = 3
$totals.Cells.Item(5, 2).Value = "2021-Q2"
$totals.Cells.Item(5, 3).Value = 2
$totals.Cells.Item(5, 4).Value = 2.76

$totals.Cells.Item(4, 1).Value = 2
$totals.Cells.Item(4, 2).Value = "2021-Q3"
$totals.Cells.Item(4, 3).Value = 2
$totals.Cells.Item(4, 4).Value = 3.17

$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(3, 2).Value = "2021-Q4"
$totals.Cells.Item(3, 3).Value = 2
$totals.Cells.Item(3, 4).Value = 2.81

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q1"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 2.82
